$wb = $excel.ActiveWorkbook

# --- GET_equipment (sheet1): selection moves from H8 to H2. Cell content unchanged. ---
$wsEquipment = $wb.Worksheets.Item("GET_equipment")
$wsEquipment.Range("H2").Select() | Out-Null

# --- GET_last_login (sheet2): fill in the schemaValidationFile cell (H2), previously blank. ---
$wsLastLogin = $wb.Worksheets.Item("GET_last_login")
$wsLastLogin.Range("H2").Value = "GET_last_login.json"

# --- GET_equipment_session (sheet3): fill in schemaValidationFile (H2, previously blank),
#     set page orientation to portrait (adds pageSetup), and it stops being the active/tabSelected sheet. ---
$wsSession = $wb.Worksheets.Item("GET_equipment_session")
$wsSession.Range("H2").Value = "GET_equipment_session.json"
$wsSession.PageSetup.Orientation = 1

# --- GET_equipment_summary (sheet4): Description (C2) & schemaValidationFile (H2) corrected,
#     page orientation set to portrait, and this sheet becomes the active/selected one. ---
$wsSummary = $wb.Worksheets.Item("GET_equipment_summary")
$wsSummary.Range("C2").Value = "Get equipment session"
$wsSummary.Range("H2").Value = "GET_equipment_summary.json"
$wsSummary.PageSetup.Orientation = 1
$wsSummary.Range("J21").Select() | Out-Null
